$d = $word.ActiveDocument
$d.Content.Find.Execute("Utente_Anonimo” invece", $true, $false, $false, $false, $false, $true, 1, $false, "Utente_Autenticato” invece", 2)
